# Update cryptos list figures (price + 1h volume change) pulled from the latest run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.083.06'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '1.665.10'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.76%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.99'
$ws.Range('E5').Value = '  -4.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5154'
$ws.Range('E6').Value = '  -5.64%  '
$ws.Range('E7').Value = '  -0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2630'
$ws.Range('E8').Value = '  -3.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06210'
$ws.Range('E9').Value = '  -3.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.06'
$ws.Range('E10').Value = '  -4.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07506'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('D12').Value = '1.667.09'
$ws.Range('E12').Value = '  -1.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.420'
$ws.Range('E13').Value = '  -2.55%  '
$ws.Range('E14').Value = '  -4.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000007923'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.53'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').Value = '26.102.87'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.771'
$ws.Range('E19').Value = '  -3.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.37'
$ws.Range('E20').Value = '  -5.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '186.06'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.124'
$ws.Range('E22').Value = '  -1.97%  '
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '146.75'
$ws.Range('E24').Value = '  -1.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1238'
$ws.Range('E25').Value = '  -5.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.518'
$ws.Range('E26').Value = '  -4.48%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06297'
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.345'
$ws.Range('E29').Value = '  -4.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.269'
$ws.Range('E30').Value = '  -4.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.465'
$ws.Range('E31').Value = '  -3.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.416'
$ws.Range('E32').Value = '  -4.81%  '
$ws.Range('E33').Value = '  -4.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9934'
$ws.Range('E34').Value = '  -4.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.406'
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6009'
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.701'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.085'
$ws.Range('E38').Value = '  -3.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01606'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').Value = '1.078.14'
$ws.Range('E40').Value = '  -3.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8564'
$ws.Range('E41').Value = '  -2.76%  '
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.92'
$ws.Range('E43').Value = '  -2.46%  '
$ws.Range('D44').Value = '1.812.85'
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000109'
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.88'
$ws.Range('E46').Value = '  -2.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.003'
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05249'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.913'
$ws.Range('E49').Value = '  -3.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4247'
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.890'
$ws.Range('E51').Value = '  -2.47%  '
